$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# 1. Fix the local raw-data path: projects moved from raw/API to raw/ecotaxa.
#    Project_localpath (column C) holds this same string in every data row.
$ws.Range("C2:C47").Value = "~/GIT/PSSdb/raw/ecotaxa/IFCB"

# 2. The author manually re-sized column C (Project_localpath) narrower and let
#    the far-right columns (D:AS) fall back to the default column width.
$ws.Columns.Item(3).ColumnWidth = 24
$ws.Range("D1:AS1").EntireColumn.ColumnWidth = 8

# 3. Update the on-screen selection/scroll state left behind after editing,
#    matching where the author's cursor ended up (column AA area, cell AM23).
$ws.Activate()
$ws.Range("AM23").Select() | Out-Null
